{"js": "// Remove the \"Ver no Jupiter ...\" paragraph, the \"\u00a9 2020 ...\" paragraph\n// that follows it, and the blank paragraph that precedes it (the trio of\n// paragraphs that used to sit right after the \"Requisitos\" section's last\n// entry, just before the trailing blank paragraph / page break).\n\nconst body = context.document.body;\n\n// Locate the \"Ver no Jupiter\" paragraph by its literal text so the script\n// does not depend on fixed paragraph indices.\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n  const paragraphs = hit.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const verParagraph = paragraphs.items[0];\n  const blankParagraph = verParagraph.getPrevious();\n  const copyrightParagraph = verParagraph.getNext();\n  await context.sync();\n\n  // Delete the copyright paragraph and the \"Ver no Jupiter\" paragraph and\n  // the preceding blank paragraph.\n  copyrightParagraph.delete();\n  verParagraph.delete();\n  blankParagraph.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter ...\" paragraph, the \"\u00a9 2020 ...\" paragraph\n# that follows it, and the blank paragraph that precedes it (the trio of\n# paragraphs that used to sit right after the \"Requisitos\" section's last\n# entry, just before the trailing blank paragraph / page break).\n\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\n\nif ($found) {\n    # Pin down the Paragraph object that contains the found text (rather\n    # than relying on a fixed paragraph index) so we can reliably reach its\n    # neighbouring paragraphs.\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {\n            $target = $p\n            break\n        }\n    }\n\n    if ($target -ne $null) {\n        $prevPara = $target.Previous()\n        $nextPara = $target.Next()\n\n        # One contiguous range spanning [blank paragraph] .. [\"Ver no\n        # Jupiter...\"] .. [\"\u00a9 2020...\"] \u2014 delete it in a single shot so the\n        # paragraph marks disappear along with the text.\n        $delRange = $d.Range($prevPara.Range.Start, $nextPara.Range.End)\n        $delRange.Delete()\n    }\n}\n"}
